$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("TestUser_1ASltkyXYuI", "automation_test+1611390700+afMLZuyHBy@gmail.com", "TestPassw0rd@123!`$wAFQkDdmos"),
    @("TestUser_1mevqxjIwaa", "automation_test+1611426616+mzWnKOpXgz@gmail.com", "TestPassw0rd@123!`$bSYYzpxNUl"),
    @("TestUser_1EpFdZwhfiQ", "automation_test+1611426867+bLXuuqccZE@gmail.com", "TestPassw0rd@123!`$fIRELWtCja")
)

$row = 4
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
